# Adds a new last column (BB, column 54) to the forecast sheet:
#  - BB1 gets the next period date serial (45986)
#  - BB3:BB18 copy the prior last column's (BA) value for that row
#  - BB19:BB21 get newly-computed forecast values (not a straight copy)
#  - Row 2 and row 22 get no new value in BB, just an expanded row span

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCol = 54   # column BB
$oldCol = 53   # column BA

# New header date for the added forecast period
$ws.Cells.Item(1, $newCol).Value = 45986

# Match the formatting (date style with border) used by the rest of the header row
$ws.Cells.Item(1, $oldCol).Copy() | Out-Null
$ws.Cells.Item(1, $newCol).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rows 3-18: straight copy of the BA value into the new BB column
$copyRows = 3..18
foreach ($r in $copyRows) {
    $val = $ws.Cells.Item($r, $oldCol).Value2
    if ($null -ne $val) {
        $ws.Cells.Item($r, $newCol).Value = $val
    }
}

# Rows 19-21: newly computed forecast values
$ws.Cells.Item(19, $newCol).Value = 2.043309689777173
$ws.Cells.Item(20, $newCol).Value = 1.002299702378884
$ws.Cells.Item(21, $newCol).Value = 1.046422855779872
